$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 28
$ws.Cells.Item($row, 1).Value = "Record"
$ws.Cells.Item($row, 2).Value = "Balanço Geral"
$ws.Cells.Item($row, 3).Value = "Saúde"
$ws.Cells.Item($row, 4).Value = "2025-04-01T13:22"
$ws.Cells.Item($row, 5).Value = "Neutro"
$ws.Cells.Item($row, 6).Value = "Idosa mordida por cachorro teve que tomar vacina contra a raiva. *nota coberta*"
